$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-08 Sunday" "2025-06-09 Monday"

Replace-Text "756×9=" "790×2="
Replace-Text "928×9=" "997×2="
Replace-Text "358×5=" "120×5="
Replace-Text "681×4=" "728×3="
Replace-Text "431×9=" "699×4="

Replace-Text "253×3=" "430×6="
Replace-Text "904×7=" "738×4="
Replace-Text "801×7=" "532×2="
Replace-Text "447×5=" "711×5="
Replace-Text "920×2=" "659×6="

Replace-Text "322×4=" "343×8="
Replace-Text "118×2=" "728×6="
Replace-Text "910×9=" "405×4="
Replace-Text "948×3=" "641×5="
Replace-Text "922×4=" "611×3="

Replace-Text "740×2=" "927×6="
Replace-Text "255×3=" "382×2="
Replace-Text "874×9=" "350×3="
Replace-Text "298×4=" "992×7="
Replace-Text "169×6=" "798×5="

Replace-Text "642×7=" "147×6="
Replace-Text "356×3=" "243×6="
Replace-Text "190×6=" "664×9="
Replace-Text "485×9=" "163×7="
Replace-Text "450×7=" "709×3="
